$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.778.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.082.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.76"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.19%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +0.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0789"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.18%  "

$ws.Range("E11").Value = "  +2.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.390.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.773"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.54%  "

$ws.Range("E16").Value = "  +1.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.061.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.705.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.43%  "

$ws.Range("E21").Value = "  +1.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("E24").Value = "  -1.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("E27").Value = "  +5.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.83%  "

$ws.Range("E30").Value = "  -2.62%  "

$ws.Range("E31").Value = "  +2.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0633"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.50%  "

$ws.Range("E37").Value = "  -3.38%  "

$ws.Range("E38").Value = "  +0.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0979"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.32%  "

$ws.Range("E42").Value = "  +1.14%  "

$ws.Range("E43").Value = "  -2.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.462.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.85%  "

$ws.Range("E46").Value = "  +0.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.64%  "

$ws.Range("E48").Value = "  +2.13%  "

$ws.Range("E49").Value = "  +2.96%  "

$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.275.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.01%  "
